# ADD LABELING FOR MAIN INTEREST GROUP
#
# This script adds a new "MAJOR-INVOLVED" column (H) to the
# Interest-group-09-15 worksheet, populating it with a 0/1 flag for
# every data row (2-260), and corrects a mis-classified data row
# (row 50) whose ethnicity columns (D/E) had been entered as
# LATINO-1/L instead of BLACK-1/B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H
$ws.Cells.Item(1, 8).Value = "MAJOR-INVOLVED"

# Correct the data entry error in row 50 (D/E columns)
$ws.Range("D50").Value = "BLACK-1"
$ws.Range("E50").Value = "B"

# 0/1 flag values for rows 2 through 260 of the new MAJOR-INVOLVED column
$hValues = @(
    1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,1,1,
    1,1,1,1,1,1,1,1,1,1,1,0,1,1,1,1,1,1,1,0,
    0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,
    1,1,1,1,1,1,0,1,1,1,1,1,1,1,1,0,0,1,1,1,
    1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,
    1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,
    1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,
    1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,
    1,1,1,1,1,1,0,1,1,1,1,1,1,1,1,1,1,1,0,1,
    1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,
    0,1,1,1,1,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,
    1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,
    1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1
)

for ($i = 0; $i -lt $hValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $hValues[$i]
}

# Match the author's final view state: zoomed to 85%, scrolled near the
# bottom of the sheet, with the last populated cell (H260) selected.
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 238
$win.ScrollColumn = 1
$ws.Range("H260").Select()
